# read/write transaksi - bill oke!
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing "kode transaksi" values (column B) were renumbered: every
# transaction code from row 9 down to row 32 shifts up by 4
# (e.g. #ID09001 -> #ID09005, #ID26006 -> #ID26010, ... #ID26022 -> #ID26026).
$bUpdates = @{
    9  = "#ID09005"
    10 = "#ID09005"
    11 = "#ID09006"
    12 = "#ID09006"
    13 = "#ID09007"
    14 = "#ID09008"
    15 = "#ID09009"
    16 = "#ID26010"
    17 = "#ID26011"
    18 = "#ID26012"
    19 = "#ID26013"
    20 = "#ID26014"
    21 = "#ID26015"
    22 = "#ID26016"
    23 = "#ID26017"
    24 = "#ID26018"
    25 = "#ID26019"
    26 = "#ID26020"
    27 = "#ID26021"
    28 = "#ID26022"
    29 = "#ID26023"
    30 = "#ID26024"
    31 = "#ID26025"
    32 = "#ID26026"
}

foreach ($row in $bUpdates.Keys) {
    $ws.Cells.Item($row, 2).Value = $bUpdates[$row]
}

# New transaction #ID30027 recorded for "Arasseo" (rows 33-34)
$ws.Cells.Item(33, 1).Value = 43281.922037662036
$ws.Cells.Item(33, 2).Value = "#ID30027"
$ws.Cells.Item(33, 3).Value = "Arasseo"
$ws.Cells.Item(33, 4).Value = 491150.00000000006

$ws.Cells.Item(34, 1).Value = 43281.922037662036
$ws.Cells.Item(34, 2).Value = "#ID30027"
$ws.Cells.Item(34, 3).Value = "PHD"
$ws.Cells.Item(34, 4).Value = 91300.0

$ws.Range("A1").Select()
